$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6: fix merged column issue - A6 was missing ("present")
$ws.Range("A6").Value = "present"

# Row 14: was a duplicate of row7 (extant->present); correct it to unknown -> uncertain
$ws.Range("A14").Value = "unknown"
$ws.Range("B14").Value = "uncertain"

# Row 15: was "unknown" -> "uncertain"; correct to "eradicated" -> "absent"
$ws.Range("A15").Value = "eradicated"
$ws.Range("B15").Value = "absent"

# New rows 16 and 17 with additional terms
$ws.Range("A16").Value = "not known"
$ws.Range("B16").Value = "uncertain"

$ws.Range("A17").Value = "data deficient"
$ws.Range("B17").Value = "uncertain"
